$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Turn D88:D93 into a shared formula group (t="shared" ref="D88:D93" si="...")
$ws.Range("D88:D93").Formula = "=C88/(24*60)"

# --- Row 94: 2025-10-17, 四方坪站 ---
$ws.Range("A94").Value = 45947
$ws.Range("B94").Value = "四方坪站"
$ws.Range("C94").Formula = "=15588/126"
$ws.Range("D94").Formula = "=C94/(24*60)"
$ws.Range("E94").Formula = "=8361.11/126"
$ws.Range("F94").Formula = "=2890.07/126"
$ws.Range("G94").Formula = "=8361.11/(15588/60)"
$ws.Range("H94").Formula = "=373/126"

# --- Row 95: 2025-10-17, 高岭站 ---
$ws.Range("A95").Value = 45947
$ws.Range("B95").Value = "高岭站"
$ws.Range("C95").Formula = "=5103/36"
$ws.Range("D95").Formula = "=C95/(24*60)"
$ws.Range("E95").Formula = "=3470.37/36"
$ws.Range("F95").Formula = "=974.57/36"
$ws.Range("G95").Formula = "=3470.37/(5103/60)"
$ws.Range("H95").Formula = "=131/36"

# Update the tracked selection to match the authored state
$ws.Range("I97").Select()
